# Apply weekly crime-data refresh to the 94th Precinct CompStat report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# A8 holds "Volume 30   Number  25" as rich text; only the trailing issue
# number run ("25" -> "26") changes.
$ws.Range("A8").Characters(21, 2).Text = "26"

# C9 holds "Report Covering the Week  6/19/2023  Through  6/25/2023" as
# rich text; the two date runs move forward by one week.
$ws.Range("C9").Characters(27, 9).Text = "6/26/2023"
$ws.Range("C9").Characters(47, 9).Text = "7/2/2023"

# --- Row 16 (Robbery) ------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -68.75
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = -21.538461538461
$ws.Range("L16").Value = 27.5
$ws.Range("M16").Value = -23.880597014925
$ws.Range("N16").Value = -83.495145631068

# --- Row 17 (Fel. Assault) --------------------------------------------------
# C17 switches from the "N/A" text placeholder to a real number.
$ws.Range("C17").Value = 5
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 49
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = -30
$ws.Range("L17").Value = 25.641025641025
$ws.Range("M17").Value = 25.641025641025
$ws.Range("N17").Value = -59.836065573770

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -23.529411764705
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = -14.953271028037
$ws.Range("N18").Value = -83.778966131907

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 6.25
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 13.461538461538
$ws.Range("I19").Value = 333
$ws.Range("J19").Value = 272
$ws.Range("K19").Value = 22.426470588235
$ws.Range("L19").Value = 75.263157894736
$ws.Range("M19").Value = 144.852941176471
$ws.Range("N19").Value = 82.967032967033

# --- Row 20 (G.L.A.) --------------------------------------------------------
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -17.647058823529
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 73
$ws.Range("K20").Value = 16.438356164383
$ws.Range("L20").Value = 66.666666666666
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -81.277533039647

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 3.125
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -11.926605504587
$ws.Range("I21").Value = 612
$ws.Range("J21").Value = 603
$ws.Range("K21").Value = 1.492537313432
$ws.Range("L21").Value = 42.990654205607
$ws.Range("M21").Value = 46.411483253588
$ws.Range("N21").Value = -62.568807339449

# --- Row 22 (Transit) --------------------------------------------------------
# C22, D22 and F22 switch from "N/A" text to real numbers; E22 switches from
# the "***.*" text placeholder to a real percentage number.
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 7
$ws.Range("L22").Value = 75
$ws.Range("M22").Value = 0

# --- Row 24 (Petit Larceny) -------------------------------------------------
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 13.043478260869
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = -11.111111111111
$ws.Range("I24").Value = 467
$ws.Range("J24").Value = 507
$ws.Range("K24").Value = -7.889546351084
$ws.Range("L24").Value = 21.932114882506
$ws.Range("M24").Value = 74.253731343283

# --- Row 25 (Misd. Assault) -------------------------------------------------
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 113
$ws.Range("J25").Value = 144
$ws.Range("K25").Value = -21.527777777777
$ws.Range("L25").Value = 15.306122448979
$ws.Range("M25").Value = 1.801801801801

# --- Row 27 (Other Sex Crimes) ----------------------------------------------
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = 76.923076923076
$ws.Range("L27").Value = 64.285714285714

# --- Row 30 (Hate Crimes) ----------------------------------------------------
# C30 and F30 switch from "N/A" text to real numbers.
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 2
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("I30").Value = 6
$ws.Range("L30").Value = 50
